$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")
$summary = $wb.Worksheets.Item("Summary")

# Updated MyForecast (D), Inventory Coverage (H), Stockout Risk (I),
# Reorder Urgency (J), Seasonality Index (L) values reflecting the
# newly added 4-week low sales check.

$ws.Range("L2").Value = 1.19

$ws.Range("D3").Value = 97
$ws.Range("H3").Value = 10.6
$ws.Range("L3").Value = 1.15

$ws.Range("D4").Value = 97
$ws.Range("H4").Value = 9.6
$ws.Range("L4").Value = 1.12

$ws.Range("D5").Value = 98
$ws.Range("H5").Value = 8.51
$ws.Range("L5").Value = 0.98

$ws.Range("D6").Value = 99
$ws.Range("H6").Value = 7.43
$ws.Range("L6").Value = 1.16

$ws.Range("D7").Value = 99
$ws.Range("H7").Value = 6.43
$ws.Range("L7").Value = 0.85

$ws.Range("D8").Value = 101
$ws.Range("H8").Value = 5.33
$ws.Range("L8").Value = 1.18

$ws.Range("D9").Value = 101
$ws.Range("H9").Value = 4.33
$ws.Range("L9").Value = 0.97

$ws.Range("D10").Value = 102
$ws.Range("H10").Value = 3.29
$ws.Range("L10").Value = 1.08

$ws.Range("D11").Value = 103
$ws.Range("H11").Value = 2.27
$ws.Range("L11").Value = 1.1

$ws.Range("D12").Value = 103
$ws.Range("H12").Value = 1.27
$ws.Range("L12").Value = 0.9399999999999999

$ws.Range("D13").Value = 105
$ws.Range("H13").Value = 0.27
$ws.Range("I13").Value = "High"
$ws.Range("J13").Value = "Urgent"
$ws.Range("L13").Value = 0.87

$ws.Range("D14").Value = 106
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = "High"
$ws.Range("L14").Value = 0.98

$ws.Range("D15").Value = 107
$ws.Range("L15").Value = 1.14

$ws.Range("D16").Value = 107
$ws.Range("L16").Value = 1.06

$ws.Range("D17").Value = 107
$ws.Range("L17").Value = 1.07

# Recalculated summary statistics on the Summary sheet.
# NumberFormat is forced to text ("@") first so the values are stored
# as text, matching the existing inline-string cells on this sheet.
$summary.Range("B9").NumberFormat = "@"
$summary.Range("B9").Value = "1628"

$summary.Range("B10").NumberFormat = "@"
$summary.Range("B10").Value = "788"

$summary.Range("B11").NumberFormat = "@"
$summary.Range("B11").Value = "388"

$summary.Range("B12").NumberFormat = "@"
$summary.Range("B12").Value = "107"

$summary.Range("B14").NumberFormat = "@"
$summary.Range("B14").Value = "96"
